$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4, column A: the question text moves to reuse the shared string
# that used to be "I tend to feel a moral obligation to meet higher
# standards than I expect from others."
$ws.Range("A4").Value = "`"I tend to feel a moral obligation to meet higher standards than I expect from others.`""

# Clear out the contents of rows 5-10 (A:C) -- the extra quiz questions
# were removed, leaving the rows present but empty.
$ws.Range("A5:C10").ClearContents()

# Update the selection to match the new active cell / selected range.
$null = $ws.Range("A5:C10").Select()
